# Going through the dataset, updating.
# Fills in the previously-unknown fields for row 2 (the Lianhuaqingwen
# capsule RCT reference) and adds a new "Other found locations" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column I.
$ws.Range("I1").Value = "Other found locations"

# Row 2: replace the placeholder "unknown" values with the real,
# looked-up bibliographic data.
$ws.Range("C2").Value = "Efficacy and safety of Lianhuaqingwen capsules, a repurposed Chinese herb, in patients with coronavirus disease 2019: A multicenter, prospective, randomized controlled trial"

$ws.Range("D2").Value = "Background`nCoronavirus disease 2019 (Covid-19) has resulted in a global outbreak.`n Few existing targeted medications are available.`n Lianhuaqingwen (LH) capsule, a repurposed marketed Chinese herb product, has been proven effective for influenza.`nPurpose`nTo determine the safety and efficacy of LH capsule in patients with Covid-19.`nMethods`nWe did a prospective multicenter open-label randomized controlled trial on LH capsule in confirmed cases with Covid-19. Patients were randomized to receive usual treatment alone or in combination with LH capsules (4 capsules, thrice daily) for 14 days.`n The primary endpoint was the rate of symptom (fever, fatigue, coughing) recovery.`nResults`nWe included 284 patients (142 each in treatment and control group) in the full-analysis set.`n The recovery rate was significantly higher in treatment group as compared with control group (91.5% vs.`n 82.4%, p = 0.022).`n The median time to symptom recovery was markedly shorter in treatment group (median: 7 vs.`n 10 days, p < 0.001).`n Time to recovery of fever (2 vs.`n 3 days), fatigue (3 vs.`n 6 days) and coughing (7 vs.`n 10 days) was also significantly shorter in treatment group (all p < 0.001).`n The rate of improvement in chest computed tomographic manifestations (83.8% vs.`n 64.1%, p < 0.001) and clinical cure (78.9% vs.`n 66.2%, p = 0.017) was also higher in treatment group.`n However, both groups did not differ in the rate of conversion to severe cases or viral assay findings (both p > 0.05).`n No serious adverse events were reported.`nConclusion`nIn light of the safety and effectiveness profiles, LH capsules could be considered to ameliorate clinical symptoms of Covid-19.`n"

$ws.Range("E2").Value = "[Ke%Hu%NULL%1, Wei-jie%Guan%NULL%0, Ying%Bi%NULL%1, Wei%Zhang%NULL%0, Lanjuan%Li%NULL%1, Boli%Zhang%NULL%1, Qingquan%Liu%NULL%0, Yuanlin%Song%NULL%1, Xingwang%Li%NULL%1, Zhongping%Duan%NULL%1, Qingshan%Zheng%NULL%1, Zifeng%Yang%NULL%1, Jingyi%Liang%NULL%1, Mingfeng%Han%NULL%1, Lianguo%Ruan%NULL%1, Chaomin%Wu%NULL%1, Yunting%Zhang%NULL%1, Zhen-hua%Jia%NULL%1, Nan-shan%Zhong%NULL%1]"

$ws.Range("F2").Value = "PMC7229744"
$ws.Range("G2").Value = "PMC"

# H2 holds a date-shaped string ("2020-05-08"). Writing it straight into
# Range.Value would make Excel auto-convert it to a date serial number,
# which isn't what the source file stores (H2 stays a shared-string cell
# with no number format). Stage it on a scratch cell formatted as text,
# then copy that value+format over so H2 keeps a plain text value.
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "2020-05-08"
$ws.Range("Z1").Copy()
$ws.Range("H2").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

# New value for the new column I in row 2.
$ws.Range("I2").Value = "_elsevier"

# The long multi-line abstract in D2 triggers an automatic row-height
# bump; restore row 2 to its natural (non-custom) height.
$ws.Rows(2).EntireRow.AutoFit()
